# Task_list_sample.xlsx -- "add rerun func, and refine table css show all"
#
# The sheet originally listed 4 rows (3 cases). This edit:
#   1. Appends two brand-new case numbers (rows 5 & 6).
#   2. "Re-runs" the existing "250400031HZH" case several more times so the
#      table is filled all the way down to row 20 (rows 7-20), matching the
#      assistant/engineer pairing used throughout the sheet.
#   3. Leaves the cursor/selection on A11, where editing stopped.
#
# Net effect on the OOXML: dimension grows from A1:C4 to A1:C20, two new
# shared strings are interned ("250404031HZH", "250405031HZH"), and the
# sheetView selection moves to A11.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -- Row 5: first new case -------------------------------------------------
$ws.Range("A5").Value = "250404031HZH"
$ws.Range("B5").Value = "Sylvia Wang"
$ws.Range("C5").Value = "Jack Chen"

# -- Row 6: second new case -------------------------------------------------
$ws.Range("A6").Value = "250405031HZH"
$ws.Range("B6").Value = "Sylvia Wang"
$ws.Range("C6").Value = "Jack Chen"

# -- Rows 7-20: rerun the "250400031HZH" case to fill out the table --------
for ($r = 7; $r -le 20; $r++) {
    $ws.Range("A$r").Value = "250400031HZH"
    $ws.Range("B$r").Value = "Sylvia Wang"
    $ws.Range("C$r").Value = "Jack Chen"
}

# Leave the selection where the author's editing session ended.
$ws.Range("A11").Select()
